# Applies the edit described by the diff:
#  - refresh several "want to go" (F) counters and one cover image URL (I2) on existing rows
#  - insert one brand-new row (event) at row 19 on both the "展览" and "全部类型" sheets,
#    shifting the remaining rows down by one and keeping column A sequential

$wb = $excel.ActiveWorkbook

# ---- Worksheet #1 ----
$ws1 = $wb.Worksheets.Item(1)

$ws1Updates = @{
    2 = @{ F = 1575; I = "//i1.hdslb.com/bfs/openplatform/202406/1UVGJ3G01718620439056.jpeg" }
    3 = @{ F = 8917 }
    4 = @{ F = 100 }
    5 = @{ F = 500 }
    6 = @{ F = 671 }
    7 = @{ F = 334 }
    8 = @{ F = 165 }
    9 = @{ F = 40 }
    10 = @{ F = 56 }
    11 = @{ F = 3778 }
    13 = @{ F = 375 }
    15 = @{ F = 4122 }
    16 = @{ F = 5 }
    18 = @{ F = 1135 }
}
foreach ($row in $ws1Updates.Keys) {
    foreach ($col in $ws1Updates[$row].Keys) {
        $ws1.Range("$col$row").Value2 = $ws1Updates[$row][$col]
    }
}

# Insert one new row at position 19 (shifts all following rows down by one).
$ws1.Rows.Item(19).Insert()
# Clone column-A formatting (bold, centered, bordered) from the row above into the new row.
$ws1.Range("A18").Copy()
$ws1.Range("A19").PasteSpecial(-4122)
$ws1.Application.CutCopyMode = $false

$ws1NewRows = @(
    @(18, "2024-07-27", "合肥·灵能百分百ONLY2.0", "铜陵北路金邦国际大厦一楼 格律诗婚礼艺术中心(新站店)", "2024.07.27 10:00-07.27 17:00", 2, 75, "https://show.bilibili.com/platform/detail.html?id=87497", "//i1.hdslb.com/bfs/openplatform/202406/3Jycwu1U1717858639976.jpeg"),
    @(19, "2024-07-27", "安徽·MAX特摄only展", "桐城路127号合作经济广场3号楼23层 赤阑桥艺术空间", "2024.07.27 09:30-07.27 18:00", 328, 50, "https://show.bilibili.com/platform/detail.html?id=83684", "//i0.hdslb.com/bfs/openplatform/202405/qBnW1VeB1715423018997.jpeg"),
    @(20, "2024-07-28", "合肥·咒术回战only", "清河路19号 依立腾工业园区", "2024.07.28 09:30-07.28 17:30", 236, 60, "https://show.bilibili.com/platform/detail.html?id=86520", "//i2.hdslb.com/bfs/openplatform/202405/cLCM0a1e1716952386781.png"),
    @(21, "2024-07-28", "合肥·第二届TH动漫游戏嘉年华", "田埠西路199号 吉祥如意宴会楼蜀山店", "2024.07.28 09:30-07.28 17:00", 8, 55, "https://show.bilibili.com/platform/detail.html?id=87447", "//i0.hdslb.com/bfs/openplatform/202406/jHqfdzLQ1718091324240.png"),
    @(22, "2024-08-03", "合肥·第七届环形宇宙动漫游戏嘉年华", "南京路与庐州大道交汇处 合肥滨湖国际会展中心", "2024.08.03 09:30-08.04 17:00", 2583, 49, "https://show.bilibili.com/platform/detail.html?id=84767", "//i2.hdslb.com/bfs/openplatform/202404/nBGuQecO1713856894035.jpeg"),
    @(23, "2024-08-17", "合肥·银魂主题派对only2.0", "长江东路1137号圣大国际商贸中心2-301室 梦田音乐LiveHouse(合肥店)", "2024.08.17 13:00-08.17 18:00", 99, 128, "https://show.bilibili.com/platform/detail.html?id=87173", "//i2.hdslb.com/bfs/openplatform/202406/aSc8SoTl1718078234193.png"),
)
$r = 19
foreach ($rowData in $ws1NewRows) {
    $ws1.Range("A" + $r).Value2 = $rowData[0]
    # Column B holds plain-text dates in the source file (not real Excel dates),
    # so force Text format first to stop auto-conversion into a date serial number.
    $ws1.Range("B" + $r).NumberFormat = "@"
    $ws1.Range("B" + $r).Value2 = $rowData[1]
    $ws1.Range("C" + $r).Value2 = $rowData[2]
    $ws1.Range("D" + $r).Value2 = $rowData[3]
    $ws1.Range("E" + $r).Value2 = $rowData[4]
    $ws1.Range("F" + $r).Value2 = $rowData[5]
    $ws1.Range("G" + $r).Value2 = $rowData[6]
    $ws1.Range("H" + $r).Value2 = $rowData[7]
    $ws1.Range("I" + $r).Value2 = $rowData[8]
    $r = $r + 1
}

# ---- Worksheet #4 ----
$ws4 = $wb.Worksheets.Item(4)

$ws4Updates = @{
    2 = @{ F = 1575; I = "//i1.hdslb.com/bfs/openplatform/202406/1UVGJ3G01718620439056.jpeg" }
    3 = @{ F = 8917 }
    4 = @{ F = 100 }
    5 = @{ F = 500 }
    6 = @{ F = 671 }
    7 = @{ F = 334 }
    8 = @{ F = 165 }
    9 = @{ F = 40 }
    10 = @{ F = 56 }
    11 = @{ F = 3778 }
    13 = @{ F = 375 }
    15 = @{ F = 4122 }
    16 = @{ F = 5 }
    18 = @{ F = 1135 }
}
foreach ($row in $ws4Updates.Keys) {
    foreach ($col in $ws4Updates[$row].Keys) {
        $ws4.Range("$col$row").Value2 = $ws4Updates[$row][$col]
    }
}

# Insert one new row at position 19 (shifts all following rows down by one).
$ws4.Rows.Item(19).Insert()
# Clone column-A formatting (bold, centered, bordered) from the row above into the new row.
$ws4.Range("A18").Copy()
$ws4.Range("A19").PasteSpecial(-4122)
$ws4.Application.CutCopyMode = $false

$ws4NewRows = @(
    @(18, "2024-07-27", "合肥·灵能百分百ONLY2.0", "铜陵北路金邦国际大厦一楼 格律诗婚礼艺术中心(新站店)", "2024.07.27 10:00-07.27 17:00", 2, 75, "https://show.bilibili.com/platform/detail.html?id=87497", "//i1.hdslb.com/bfs/openplatform/202406/3Jycwu1U1717858639976.jpeg"),
    @(19, "2024-07-27", "安徽·MAX特摄only展", "桐城路127号合作经济广场3号楼23层 赤阑桥艺术空间", "2024.07.27 09:30-07.27 18:00", 328, 50, "https://show.bilibili.com/platform/detail.html?id=83684", "//i0.hdslb.com/bfs/openplatform/202405/qBnW1VeB1715423018997.jpeg"),
    @(20, "2024-07-28", "合肥·咒术回战only", "清河路19号 依立腾工业园区", "2024.07.28 09:30-07.28 17:30", 236, 60, "https://show.bilibili.com/platform/detail.html?id=86520", "//i2.hdslb.com/bfs/openplatform/202405/cLCM0a1e1716952386781.png"),
    @(21, "2024-07-28", "合肥·第二届TH动漫游戏嘉年华", "田埠西路199号 吉祥如意宴会楼蜀山店", "2024.07.28 09:30-07.28 17:00", 8, 55, "https://show.bilibili.com/platform/detail.html?id=87447", "//i0.hdslb.com/bfs/openplatform/202406/jHqfdzLQ1718091324240.png"),
    @(22, "2024-08-03", "合肥·第七届环形宇宙动漫游戏嘉年华", "南京路与庐州大道交汇处 合肥滨湖国际会展中心", "2024.08.03 09:30-08.04 17:00", 2583, 49, "https://show.bilibili.com/platform/detail.html?id=84767", "//i2.hdslb.com/bfs/openplatform/202404/nBGuQecO1713856894035.jpeg"),
    @(23, "2024-08-03", "合肥·首届包河留声机音乐节—《菊次郎的夏天》久石让钢琴曲梦幻之旅演奏会", "徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院", "2024.08.03 19:30-08.03 21:00", 37, 80, "https://show.bilibili.com/platform/detail.html?id=83556", "//i1.hdslb.com/bfs/openplatform/202403/4nwOTVDu1711695345941.jpeg"),
    @(24, "2024-08-17", "合肥·银魂主题派对only2.0", "长江东路1137号圣大国际商贸中心2-301室 梦田音乐LiveHouse(合肥店)", "2024.08.17 13:00-08.17 18:00", 99, 128, "https://show.bilibili.com/platform/detail.html?id=87173", "//i2.hdslb.com/bfs/openplatform/202406/aSc8SoTl1718078234193.png"),
)
$r = 19
foreach ($rowData in $ws4NewRows) {
    $ws4.Range("A" + $r).Value2 = $rowData[0]
    # Column B holds plain-text dates in the source file (not real Excel dates),
    # so force Text format first to stop auto-conversion into a date serial number.
    $ws4.Range("B" + $r).NumberFormat = "@"
    $ws4.Range("B" + $r).Value2 = $rowData[1]
    $ws4.Range("C" + $r).Value2 = $rowData[2]
    $ws4.Range("D" + $r).Value2 = $rowData[3]
    $ws4.Range("E" + $r).Value2 = $rowData[4]
    $ws4.Range("F" + $r).Value2 = $rowData[5]
    $ws4.Range("G" + $r).Value2 = $rowData[6]
    $ws4.Range("H" + $r).Value2 = $rowData[7]
    $ws4.Range("I" + $r).Value2 = $rowData[8]
    $r = $r + 1
}

Write-Host "Sheet1 (展览) used rows:" $ws1.UsedRange.Rows.Count
Write-Host "Sheet4 (全部类型) used rows:" $ws4.UsedRange.Rows.Count
